$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values keep their exact textual representation
# (avoids Excel auto-converting numeric-looking strings to numbers/dates,
# which would drop trailing zeros / alter formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '37.809.91'
$ws.Range('E2').Value = '  +7.05%  '
$ws.Range('D3').Value = '1.960.51'
$ws.Range('E3').Value = '  +3.97%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '251.74'
$ws.Range('E5').Value = '  +2.51%  '
$ws.Range('D6').Value = '0.697'
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '48.69'
$ws.Range('E8').Value = '  +13.67%  '
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  +8.85%  '
$ws.Range('D10').Value = '59.35'
$ws.Range('E10').Value = '  +7.65%  '
$ws.Range('D11').Value = '0.0774'
$ws.Range('E11').Value = '  +4.23%  '
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').Value = '15.90'
$ws.Range('E13').Value = '  +14.95%  '
$ws.Range('D14').Value = '0.847'
$ws.Range('E14').Value = '  +9.50%  '
$ws.Range('D15').Value = '2.228.02'
$ws.Range('E15').Value = '  +3.14%  '
$ws.Range('D16').Value = '5.22'
$ws.Range('E16').Value = '  +4.80%  '
$ws.Range('D17').Value = '1.937.93'
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '37.742.16'
$ws.Range('E18').Value = '  +6.81%  '
$ws.Range('D19').Value = '75.93'
$ws.Range('E19').Value = '  +3.57%  '
$ws.Range('D20').Value = '0.0₃0869'
$ws.Range('E20').Value = '  +5.61%  '
$ws.Range('D21').Value = '13.82'
$ws.Range('E21').Value = '  +7.91%  '
$ws.Range('D22').Value = '256.02'
$ws.Range('E22').Value = '  +4.71%  '
$ws.Range('D23').Value = '5.27'
$ws.Range('E23').Value = '  +2.27%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').Value = '2.53'
$ws.Range('E25').Value = '  -5.51%  '
$ws.Range('D26').Value = '170.34'
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('D27').Value = '2.14'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').Value = '9.01'
$ws.Range('E28').Value = '  +5.58%  '
$ws.Range('D29').Value = '19.22'
$ws.Range('E29').Value = '  +5.30%  '
$ws.Range('D30').Value = '0.130'
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('D31').Value = '4.65'
$ws.Range('E31').Value = '  +8.70%  '
$ws.Range('D32').Value = '0.0622'
$ws.Range('E32').Value = '  +4.81%  '
$ws.Range('D33').Value = '0.0926'
$ws.Range('E33').Value = '  +29.29%  '
$ws.Range('D34').Value = '4.39'
$ws.Range('E34').Value = '  +4.54%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.90'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('B36').Value = 'Gas'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D36').Value = '19.71'
$ws.Range('E36').Value = '  +43.91%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = '0.908'
$ws.Range('E38').Value = '  +6.76%  '
$ws.Range('D39').Value = '1.46'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '2.03'
$ws.Range('E40').Value = '  +4.99%  '
$ws.Range('D41').Value = '105.80'
$ws.Range('E41').Value = '  +7.96%  '
$ws.Range('D42').Value = '0.0229'
$ws.Range('E42').Value = '  +3.41%  '
$ws.Range('D43').Value = '17.67'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  +21.19%  '
$ws.Range('D45').Value = '1.13'
$ws.Range('E45').Value = '  +5.32%  '
$ws.Range('D46').Value = '1.368.04'
$ws.Range('E46').Value = '  +3.21%  '
$ws.Range('D47').Value = '2.45'
$ws.Range('E47').Value = '  +3.14%  '
$ws.Range('D48').Value = '0.0851'
$ws.Range('E48').Value = '  +5.37%  '
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('E50').Value = '  +18.02%  '
$ws.Range('D51').Value = '6.48'
$ws.Range('E51').Value = '  +3.28%  '
